# Updated Topics --> Monorepos & NPM Deploy
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOPICS")

# New row 52: "Deploying NPM " (mirrors the formatting of row 50 / hyperlink-style cells)
$ws.Range("A52").Value = "Deploying NPM "
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null

# New row 54: "Monoreps --> TurboRepo"
$ws.Range("A54").Value = "Monoreps --> TurboRepo"
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Move the active selection to match the edited workbook's cursor position
$ws.Range("M47").Select() | Out-Null
